# "removed IR sensors from SLAM"
# The "Infrared sensors" line item (row 23) is removed from Sheet1's
# equipment list; every row below it shifts up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Rows.Item(23).Delete() | Out-Null

# Restore the cursor position left behind by the edit.
$ws.Range("G18").Select() | Out-Null
